# Updated BOM for Camera System
# Adds two new line items (PIR sensor + DHT11 sensor) to the Bill of
# Materials worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFormat = """R""#,##0.00"

# ---- Row 6: PIR Sensor -------------------------------------------------
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Camera"
$ws.Range("C6").Value = "Passive Infra-Red (PIR) Sensor HC-SR501"
$ws.Range("C6").Font.Name = "Calibri"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 26
$ws.Range("F6").Value = 26
$ws.Range("E6:F6").NumberFormat = $currencyFormat

# ---- Row 7: DHT11 Sensor ------------------------------------------------
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Camera"
$ws.Range("C7").Value = "Temperature and Humidity Sensor DHT11"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 33.91
$ws.Range("F7").Value = 33.91
$ws.Range("E7:F7").NumberFormat = $currencyFormat

# Match the selection left behind by the editing session
$ws.Range("F7").Select()
